$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '58.924.83'
Set-TextValue $ws 'E2' '  -2.05%  '
Set-TextValue $ws 'D3' '2.311.96'
Set-TextValue $ws 'E3' '  -4.56%  '
Set-TextValue $ws 'D4' '1.00'
Set-TextValue $ws 'D5' '551.07'
Set-TextValue $ws 'E5' '  -0.56%  '
Set-TextValue $ws 'D6' '131.54'
Set-TextValue $ws 'E6' '  -3.97%  '
Set-TextValue $ws 'D7' '1.00'
Set-TextValue $ws 'E7' '  -0.04%  '
Set-TextValue $ws 'D8' '0.573'
Set-TextValue $ws 'E8' '  -2.80%  '
Set-TextValue $ws 'D9' '2.309.61'
Set-TextValue $ws 'E9' '  -4.59%  '
Set-TextValue $ws 'E10' '  -3.08%  '
Set-TextValue $ws 'E11' '  -2.73%  '
Set-TextValue $ws 'E12' '  +1.08%  '
Set-TextValue $ws 'D13' '0.336'
Set-TextValue $ws 'E13' '  -5.31%  '
Set-TextValue $ws 'D14' '24.13'
Set-TextValue $ws 'E14' '  -2.81%  '
Set-TextValue $ws 'D15' '2.721.65'
Set-TextValue $ws 'E15' '  -4.66%  '
Set-TextValue $ws 'D16' '58.878.18'
Set-TextValue $ws 'E16' '  -2.02%  '
Set-TextValue $ws 'E17' '  -3.01%  '
Set-TextValue $ws 'D18' '2.342.71'
Set-TextValue $ws 'E18' '  -4.07%  '
Set-TextValue $ws 'D19' '10.74'
Set-TextValue $ws 'E19' '  -4.57%  '
Set-TextValue $ws 'E20' '  -3.73%  '
Set-TextValue $ws 'D21' '316.18'
Set-TextValue $ws 'E21' '  -3.47%  '
Set-TextValue $ws 'E22' '  -4.27%  '
Set-TextValue $ws 'E23' '  +0.14%  '
Set-TextValue $ws 'D24' '63.29'
Set-TextValue $ws 'E24' '  -2.85%  '
Set-TextValue $ws 'D25' '0.173'
Set-TextValue $ws 'E25' '  -2.77%  '
Set-TextValue $ws 'D26' '1.00'
Set-TextValue $ws 'E26' '  -0.09%  '
Set-TextValue $ws 'D27' '8.08'
Set-TextValue $ws 'E27' '  -6.59%  '
Set-TextValue $ws 'E28' '  -7.27%  '
Set-TextValue $ws 'E29' '  +0.12%  '
Set-TextValue $ws 'D30' '170.32'
Set-TextValue $ws 'E30' '  +0.04%  '
Set-TextValue $ws 'D31' '0.0₃0733'
Set-TextValue $ws 'E31' '  -5.08%  '
Set-TextValue $ws 'E32' '  +3.30%  '
Set-TextValue $ws 'E33' '  -4.47%  '
Set-TextValue $ws 'E34' '  -4.25%  '
Set-TextValue $ws 'E35' '  +0.02%  '
Set-TextValue $ws 'E36' '  -3.65%  '
Set-TextValue $ws 'E37' '  +0.00%  '
Set-TextValue $ws 'E38' '  -6.45%  '
Set-TextValue $ws 'D39' '4.00'
Set-TextValue $ws 'E39' '  -5.35%  '
Set-TextValue $ws 'D40' '38.18'
Set-TextValue $ws 'E40' '  -1.91%  '
Set-TextValue $ws 'E41' '  -4.82%  '
Set-TextValue $ws 'D42' '305.03'
Set-TextValue $ws 'E42' '  -6.05%  '
Set-TextValue $ws 'D43' '142.58'
Set-TextValue $ws 'E43' '  -2.31%  '
Set-TextValue $ws 'D44' '3.46'
Set-TextValue $ws 'E44' '  -5.25%  '
Set-TextValue $ws 'D45' '0.0957'
Set-TextValue $ws 'E45' '  -0.67%  '
Set-TextValue $ws 'D46' '0.0504'
Set-TextValue $ws 'E46' '  -2.57%  '
Set-TextValue $ws 'D47' '18.85'
Set-TextValue $ws 'E47' '  -4.59%  '
Set-TextValue $ws 'D48' '0.561'
Set-TextValue $ws 'E49' '  -2.81%  '
Set-TextValue $ws 'D50' '16.75'
Set-TextValue $ws 'E50' '  -4.15%  '
Set-TextValue $ws 'D51' '11.02'
Set-TextValue $ws 'E51' '  -0.29%  '
